$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S mirrors column R's formatting (number format / borders / font),
# so seed it by copying R4:R34 -> S4:S34, then overwrite with the real 2022 data.
$ws.Range("R4:R34").Copy($ws.Range("S4:S34")) | Out-Null

# Header row: 2022
$ws.Range("S4").Value = 2022

# Data rows (2022 column), "-" encodes the same textual dash used elsewhere
$values = @{
    5  = 0.5
    6  = 0.2
    7  = 0.7
    8  = 0.2
    9  = "-"
    10 = 0.4
    11 = 0.5
    12 = 0.3
    13 = 0.6
    14 = 0.7
    15 = 0.4
    16 = 1.1000000000000001
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 0.4
    21 = 0.4
    22 = 0.4
    23 = 0.4
    24 = "-"
    25 = 0.7
    26 = 1
    27 = 0.4
    28 = 1.7
    29 = 0.3
    30 = 0
    31 = 0.6
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 19).Value = $values[$r]
}

# Match the author's recorded selection at save time
$ws.Range("T6").Select() | Out-Null
